$d = $word.ActiveDocument

# Insert "production, " after "natural gas " in the sentence about scaling
# total natural gas transmission/processing/distribution emissions.
$d.Content.Find.Execute(
    "data and scale the total natural gas transmission,",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "data and scale the total natural gas production, transmission,", 2)
